# Apply crypto price/volume updates per commit diff (Sun Aug 18 14:35:54 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.758.98"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "2.650.29"
$ws.Range("E3").Value = "  +2.24%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'536.58"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "'145.43"
$ws.Range("E6").Value = "  +3.78%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +1.00%  "
$ws.Range("D9").Value = "2.667.02"
$ws.Range("E9").Value = "  +2.42%  "
$ws.Range("D10").Value = "'6.75"
$ws.Range("E10").Value = "  +4.91%  "
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").Value = "3.123.54"
$ws.Range("E14").Value = "  +2.45%  "
$ws.Range("D15").Value = "59.677.97"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").Value = "'21.28"
$ws.Range("E16").Value = "  +4.21%  "
$ws.Range("D17").Value = "2.677.33"
$ws.Range("E17").Value = "  +3.28%  "
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").Value = "'345.51"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("E20").Value = "  +2.21%  "
$ws.Range("D21").Value = "'10.25"
$ws.Range("E21").Value = "  +1.60%  "
$ws.Range("D22").Value = "'6.36"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'66.80"
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").Value = "'7.30"
$ws.Range("E28").Value = "  +2.27%  "
$ws.Range("E29").Value = "  +2.61%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  +2.69%  "
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("D33").Value = "'19.07"
$ws.Range("D34").Value = "'150.24"
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("E35").Value = "  +1.90%  "
$ws.Range("E36").Value = "  +3.24%  "
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").Value = "'0.842"
$ws.Range("E38").Value = "  +1.26%  "
$ws.Range("D39").Value = "'297.54"
$ws.Range("E39").Value = "  +9.49%  "
$ws.Range("D40").Value = "'0.825"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("E41").Value = "  +2.07%  "
$ws.Range("D42").Value = "'0.998"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  +1.86%  "
$ws.Range("E44").Value = "  +5.21%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'19.42"
$ws.Range("E45").Value = "  +5.76%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "'10.73"
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("D49").Value = "1.972.24"
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("E50").Value = "  +2.55%  "
$ws.Range("D51").Value = "'18.41"
$ws.Range("E51").Value = "  +0.98%  "
